# Added periodic & upfront related scenarios
#
# The repayment-strategy value cell (B17) on the "ProductLoanInput" sheet
# changes from "RBI (India)" to "Overdue/Due Fee/Int,Principal".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Restore the scroll position / selection the author had on screen when the
# workbook was saved (best effort - cosmetic view state).
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("B17").Select()
